# "Change paths to excel files" — point the Output sheet's SUMIFS lookups at
# the external 'Sales and Inv Paste' sheet/columns instead of the stale
# hard-coded literals, and recompute the variance row from a formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: E5:G5 -> SUMIFS against 'Sales and Inv Paste' column V
$ws.Range("E5").Formula = '=SUMIFS(''Sales and Inv Paste''!$V:$V, ''Sales and Inv Paste''!$E:$E, $B$4, ''Sales and Inv Paste''!$G:$G, $E$3, ''Sales and Inv Paste''!$J:$J, E$4, ''Sales and Inv Paste''!$K:$K, $A$5)'
$ws.Range("F5").Formula = '=SUMIFS(''Sales and Inv Paste''!$V:$V, ''Sales and Inv Paste''!$E:$E, $B$4, ''Sales and Inv Paste''!$G:$G, $E$3, ''Sales and Inv Paste''!$J:$J, F$4, ''Sales and Inv Paste''!$K:$K, $A$5)'
$ws.Range("G5").Formula = '=SUMIFS(''Sales and Inv Paste''!$V:$V, ''Sales and Inv Paste''!$E:$E, $B$4, ''Sales and Inv Paste''!$G:$G, $E$3, ''Sales and Inv Paste''!$J:$J, G$4, ''Sales and Inv Paste''!$K:$K, $A$5)'

# Row 6: E6 becomes a single-space placeholder; F6/G6 stay as-is (43434)
$ws.Range("E6").Value = " "

# Row 7: E7:G7 -> SUMIFS against 'Sales and Inv Paste' column W
$ws.Range("E7").Formula = '=SUMIFS(''Sales and Inv Paste''!$W:$W, ''Sales and Inv Paste''!$E:$E, $B$4, ''Sales and Inv Paste''!$G:$G, $E$3, ''Sales and Inv Paste''!$J:$J, E$4, ''Sales and Inv Paste''!$K:$K, $A$5)'
$ws.Range("F7").Formula = '=SUMIFS(''Sales and Inv Paste''!$W:$W, ''Sales and Inv Paste''!$E:$E, $B$4, ''Sales and Inv Paste''!$G:$G, $E$3, ''Sales and Inv Paste''!$J:$J, F$4, ''Sales and Inv Paste''!$K:$K, $A$5)'
$ws.Range("G7").Formula = '=SUMIFS(''Sales and Inv Paste''!$W:$W, ''Sales and Inv Paste''!$E:$E, $B$4, ''Sales and Inv Paste''!$G:$G, $E$3, ''Sales and Inv Paste''!$J:$J, G$4, ''Sales and Inv Paste''!$K:$K, $A$5)'

# Row 8: E8:G8 -> percent-change formula, guarded against div/0
$ws.Range("E8").Formula = "=IFERROR((E5-E7)/E7,0)"
$ws.Range("F8").Formula = "=IFERROR((F5-F7)/F7,0)"
$ws.Range("G8").Formula = "=IFERROR((G5-G7)/G7,0)"
